$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on the text-like columns (Caso/OT/Comuna/etc.) for the
# rows we are about to (re)write, so Excel does not auto-coerce numeric-looking
# strings (e.g. "2098", "13") or date-looking strings (e.g. "5/24/2024") on assignment.
# NumberFormat is set one single-area range at a time: multi-area "A1,B1" range strings
# only apply to the first area in this host, so each column/row-block gets its own call.
$ws.Range("A2:A9").NumberFormat = "@"
$ws.Range("B2:B9").NumberFormat = "@"
$ws.Range("C2:C9").NumberFormat = "@"
$ws.Range("D2:D9").NumberFormat = "@"
$ws.Range("E2:E9").NumberFormat = "@"
$ws.Range("F2:F9").NumberFormat = "@"
$ws.Range("G2:G9").NumberFormat = "@"
$ws.Range("H2:H9").NumberFormat = "@"
$ws.Range("J2:J9").NumberFormat = "@"
$ws.Range("K2:K9").NumberFormat = "@"
$ws.Range("L2:L9").NumberFormat = "@"
$ws.Range("O2:O9").NumberFormat = "@"
$ws.Range("P2:P9").NumberFormat = "@"
$ws.Range("A61:A65").NumberFormat = "@"
$ws.Range("B61:B65").NumberFormat = "@"
$ws.Range("C61:C65").NumberFormat = "@"
$ws.Range("D61:D65").NumberFormat = "@"
$ws.Range("E61:E65").NumberFormat = "@"
$ws.Range("F61:F65").NumberFormat = "@"
$ws.Range("G61:G65").NumberFormat = "@"
$ws.Range("H61:H65").NumberFormat = "@"
$ws.Range("J61:J65").NumberFormat = "@"
$ws.Range("K61:K65").NumberFormat = "@"
$ws.Range("L61:L65").NumberFormat = "@"
$ws.Range("O61:O65").NumberFormat = "@"
$ws.Range("P61:P65").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = '5589'
$ws.Range("B2").Value = '12/31/2023'
$ws.Range("C2").Value = 'ARCOS 1520'
$ws.Range("D2").Value = '13'
$ws.Range("E2").Value = '799540526'
$ws.Range("F2").Value = 'NEW'
$ws.Range("G2").Value = 'Pendiente de Traspaso PROPIO'
$ws.Range("H2").Value = 'Picada'
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = 'Pasante'
$ws.Range("O2").Value = 'Colegiales'
$ws.Range("P2").Value = 'Capital Norte'
$ws.Range("I2").Value = 0
$ws.Range("M2").Value = -58.449125
$ws.Range("N2").Value = -34.565958

# Row 3
$ws.Range("A3").Value = '4862'
$ws.Range("B3").Value = '1/23/2025'
$ws.Range("C3").Value = 'ARCOS 2263'
$ws.Range("D3").Value = '13'
$ws.Range("E3").Value = '802857379'
$ws.Range("F3").Value = 'NEW'
$ws.Range("G3").Value = 'Pendiente de Traspaso PROPIO'
$ws.Range("H3").Value = 'picada'
$ws.Range("J3").Value = 'Cambio'
$ws.Range("K3").Value = 'Nodo Teco'
$ws.Range("L3").Value = 'Pasante'
$ws.Range("O3").Value = 'Saavedra'
$ws.Range("P3").Value = 'Capital Norte'
$ws.Range("I3").Value = 0
$ws.Range("M3").Value = -58.455082
$ws.Range("N3").Value = -34.558883

# Row 4
$ws.Range("A4").Value = '2098'
$ws.Range("B4").Value = '5/24/2024'
$ws.Range("C4").Value = 'AZURDUY JUANA 2449'
$ws.Range("D4").Value = '13'
$ws.Range("E4").Value = '788826017'
$ws.Range("F4").Value = 'NEW'
$ws.Range("G4").Value = 'Pendiente'
$ws.Range("H4").Value = 'Terminal con rienda'
$ws.Range("J4").Value = 'Cambio'
$ws.Range("K4").Value = 'Sin equipos'
$ws.Range("L4").Value = 'Pasante'
$ws.Range("O4").Value = 'Saavedra'
$ws.Range("P4").Value = 'Capital Norte'
$ws.Range("I4").Value = 1
$ws.Range("M4").Value = -58.467279
$ws.Range("N4").Value = -34.551117

# Row 5
$ws.Range("A5").Value = '3299'
$ws.Range("B5").Value = '9/10/2024'
$ws.Range("C5").Value = 'DIAZ COLODRERO 3309'
$ws.Range("D5").Value = '12'
$ws.Range("E5").Value = '796186684'
$ws.Range("F5").Value = 'NEW'
$ws.Range("G5").Value = 'Pendiente'
$ws.Range("H5").Value = 'Colocar columna para solicitar traspasos'
$ws.Range("J5").Value = 'Cambio'
$ws.Range("K5").Value = 'Nodo TLC'
$ws.Range("L5").Value = 'Pasante'
$ws.Range("O5").Value = 'Paternal'
$ws.Range("P5").Value = 'Capital Norte'
$ws.Range("I5").Value = 1
$ws.Range("M5").Value = -58.491722
$ws.Range("N5").Value = -34.565845

# Row 6
$ws.Range("A6").Value = '3839'
$ws.Range("B6").Value = '10/23/2024'
$ws.Range("C6").Value = 'PICO 1511'
$ws.Range("D6").Value = '13'
$ws.Range("E6").Value = '798390296'
$ws.Range("F6").Value = 'NEW'
$ws.Range("G6").Value = 'Pendiente'
$ws.Range("H6").Value = 'Poste inclinado'
$ws.Range("J6").Value = 'Aplomo'
$ws.Range("K6").Value = 'Sin equipos'
$ws.Range("L6").Value = 'Poste'
$ws.Range("O6").Value = 'Saavedra'
$ws.Range("P6").Value = 'Capital Norte'
$ws.Range("I6").Value = 1
$ws.Range("M6").Value = -58.465596
$ws.Range("N6").Value = -34.53627

# Row 7
$ws.Range("A7").Value = '801645368'
$ws.Range("B7").Value = '12/13/2024'
$ws.Range("C7").Value = 'San Blas 1809'
$ws.Range("D7").Value = '11'
$ws.Range("E7").Value = '801645368'
$ws.Range("F7").Value = 'NEW'
$ws.Range("G7").Value = 'Pendiente'
$ws.Range("H7").Value = 'Picada'
$ws.Range("J7").Value = 'Cambio'
$ws.Range("K7").Value = 'Sin equipos'
$ws.Range("L7").Value = 'Pasante'
$ws.Range("O7").Value = 'Paternal'
$ws.Range("P7").Value = 'Capital Norte'
$ws.Range("I7").Value = 0
$ws.Range("M7").Value = -58.467767
$ws.Range("N7").Value = -34.604588

# Row 8
$ws.Range("A8").Value = '4595'
$ws.Range("B8").Value = '1/15/2025'
$ws.Range("C8").Value = 'PAROISSIEN 1806'
$ws.Range("D8").Value = '13'
$ws.Range("E8").Value = '802747617'
$ws.Range("F8").Value = 'NEW'
$ws.Range("G8").Value = 'Pendiente'
$ws.Range("H8").Value = 'Aplomar'
$ws.Range("J8").Value = 'Aplomo'
$ws.Range("K8").Value = 'Sin equipos'
$ws.Range("L8").Value = 'Terminal'
$ws.Range("O8").Value = 'Saavedra'
$ws.Range("P8").Value = 'Capital Norte'
$ws.Range("I8").Value = 1
$ws.Range("M8").Value = -58.464172
$ws.Range("N8").Value = -34.543845

# Row 9
$ws.Range("A9").Value = '4662'
$ws.Range("B9").Value = '1/21/2025'
$ws.Range("C9").Value = 'ALTOLAGUIRRE 2397'
$ws.Range("D9").Value = '12'
$ws.Range("E9").Value = '802823938'
$ws.Range("F9").Value = 'NEW'
$ws.Range("G9").Value = 'Pendiente'
$ws.Range("H9").Value = 'Inclinada'
$ws.Range("J9").Value = 'Aplomo'
$ws.Range("K9").Value = 'Sin equipos'
$ws.Range("L9").Value = 'Pasante'
$ws.Range("O9").Value = 'Paternal'
$ws.Range("P9").Value = 'Capital Norte'
$ws.Range("I9").Value = 1
$ws.Range("M9").Value = -58.490766
$ws.Range("N9").Value = -34.576987

# Row 39: update Observaciones text (everything else in that row is unchanged)
$ws.Range("H39").Value = 'Picada info para cierre tambien entro como caso 6911'

# New rows 61-65 appended at the bottom of the sheet
# Row 61
$ws.Range("A61").Value = '6906'
$ws.Range("B61").Value = '8/12/2025'
$ws.Range("C61").Value = 'MOSCONI GENERAL AV. 3163'
$ws.Range("D61").Value = '12'
$ws.Range("E61").Value = '808918685'
$ws.Range("F61").Value = 'NEW'
$ws.Range("G61").Value = 'Pendiente'
$ws.Range("H61").Value = 'Picada'
$ws.Range("J61").Value = 'Cambio'
$ws.Range("K61").Value = 'Sin equipos'
$ws.Range("L61").Value = 'Pasante'
$ws.Range("O61").Value = 'Paternal'
$ws.Range("P61").Value = 'Capital Norte'
$ws.Range("I61").Value = 1
$ws.Range("M61").Value = -58.505731
$ws.Range("N61").Value = -34.588316

# Row 62
$ws.Range("A62").Value = '6910'
$ws.Range("B62").Value = '8/12/2025'
$ws.Range("C62").Value = 'RIVAS, GRAL. 2345'
$ws.Range("D62").Value = '11'
$ws.Range("E62").Value = '808918698'
$ws.Range("F62").Value = 'NEW'
$ws.Range("G62").Value = 'Pendiente'
$ws.Range("H62").Value = 'Cambiar'
$ws.Range("J62").Value = 'Cambio'
$ws.Range("K62").Value = 'Sin equipos'
$ws.Range("L62").Value = 'Pasante'
$ws.Range("O62").Value = 'Paternal'
$ws.Range("P62").Value = 'Capital Norte'
$ws.Range("I62").Value = 1
$ws.Range("M62").Value = -58.482497
$ws.Range("N62").Value = -34.598714

# Row 63
$ws.Range("A63").Value = '6928'
$ws.Range("B63").Value = '8/12/2025'
$ws.Range("C63").Value = 'ALBARELLOS AV. 3101'
$ws.Range("D63").Value = '12'
$ws.Range("E63").Value = '808918710'
$ws.Range("F63").Value = 'NEW'
$ws.Range("G63").Value = 'Pendiente'
$ws.Range("H63").Value = 'Picada'
$ws.Range("J63").Value = 'Cambio'
$ws.Range("K63").Value = 'Sin equipos'
$ws.Range("L63").Value = 'Pasante'
$ws.Range("O63").Value = 'Paternal'
$ws.Range("P63").Value = 'Capital Norte'
$ws.Range("I63").Value = 1
$ws.Range("M63").Value = -58.512623
$ws.Range("N63").Value = -34.579137

# Row 64
$ws.Range("A64").Value = '6939'
$ws.Range("B64").Value = '8/12/2025'
$ws.Range("C64").Value = 'ANDONAEGUI 1894'
$ws.Range("D64").Value = '12'
$ws.Range("E64").Value = '808918715'
$ws.Range("F64").Value = 'NEW'
$ws.Range("G64").Value = 'Pendiente'
$ws.Range("H64").Value = 'Cambiar'
$ws.Range("J64").Value = 'Cambio'
$ws.Range("K64").Value = 'Sin equipos'
$ws.Range("L64").Value = 'Pasante'
$ws.Range("O64").Value = 'Paternal'
$ws.Range("P64").Value = 'Capital Norte'
$ws.Range("I64").Value = 1
$ws.Range("M64").Value = -58.484848
$ws.Range("N64").Value = -34.581325

# Row 65
$ws.Range("A65").Value = '-550'
$ws.Range("B65").Value = '8/12/2025'
$ws.Range("C65").Value = 'Fitz roy 267'
$ws.Range("D65").Value = '15'
$ws.Range("E65").Value = '808918720'
$ws.Range("F65").Value = 'NEW'
$ws.Range("G65").Value = 'Pendiente'
$ws.Range("H65").Value = 'Aplomar columna'
$ws.Range("J65").Value = 'Cambio'
$ws.Range("K65").Value = 'Sin equipos'
$ws.Range("L65").Value = 'Pasante'
$ws.Range("O65").Value = 'Paternal'
$ws.Range("P65").Value = 'Capital Norte'
$ws.Range("I65").Value = 1
$ws.Range("M65").Value = -58.451378
$ws.Range("N65").Value = -34.596195
